$wb = $excel.ActiveWorkbook

# ---- "dog" sheet (sheet1) ----
$dog = $wb.Worksheets.Item("dog")

# Correct the dates on the existing rows (2 and 3)
$dog.Range("A2").Value = 45766
$dog.Range("A3").Value = 45767

# Add the new field-test row (row 4)
$dog.Range("A4").Value = 45772
$dog.Range("B4").Value = "PRESENCE"
$dog.Range("C4").Value = 0.4513888888888889
$dog.Range("D4").Value = 0.58333333333333337
$dog.Range("E4").Value = 23
$dog.Range("F4").Value = 2
$dog.Range("G4").Value = "Sunny, mild"
$dog.Range("H4").Value = $true
$dog.Range("I4").Value = "14 minutes 27 seconds"
$dog.Range("J4").Value = 867
$dog.Range("K4").Value = "Primary sweeps"
$dog.Range("L4").Value = "Worked uphill. Picked up odour pooling around trees uphill of target and then worked back down to target. Koda slightly hot."

# Carry the date/time number formats down into the new row
$dog.Range("A3").Copy()
$dog.Range("A4").PasteSpecial(-4122)
$dog.Range("C3:D3").Copy()
$dog.Range("C4:D4").PasteSpecial(-4122)

# Column J had a stray direct format hanging off the whole column (and
# off J2/J3) left over from earlier edits - clear it, then restore the
# bold header look on J1 only
$dog.Columns.Item(10).ClearFormats()
$dog.Range("J1").Style = "Normal"
$dog.Range("J1").Font.Bold = $true

# Make "dog" the active sheet/tab and leave the selection where the user left it
$dog.Activate()
$dog.Range("H5").Select()

# ---- "human" sheet (sheet2) ----
$human = $wb.Worksheets.Item("human")

# F2 carried the same stray direct format - clear it back to plain
$human.Range("E2").Copy()
$human.Range("F2").PasteSpecial(-4122)
